# COCOMO Estimations.xlsx - "improved the graph visualization for model and
# transaction drawers"
#
# The author's Google Drive synced folder moved from an old Windows profile
# ("C:\Users\Kan Qi\Google Drive\...") to a new one
# ("C:\Users\flyqk\Documents\Google Drive\..."), so every recorded project
# URL in column D needed that path prefix updated. The sheet also got two
# columns widened (B and C) and the active-cell selection was left on H14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the stale local file paths stored in column D (the "URL" column)
#    everywhere they appear in the used range.
$ws.UsedRange.Replace("C:\Users\Kan Qi\Google Drive", "C:\Users\flyqk\Documents\Google Drive") | Out-Null

# 2. Widen columns B (Semester) and C (PROJ) so the longer values are fully
#    visible.
$ws.Columns.Item(2).ColumnWidth = 18.3
$ws.Columns.Item(3).ColumnWidth = 41.1666666666

# 3. Leave the selection on H14, matching the saved view state.
$ws.Range("H14").Select() | Out-Null
